$wb = $excel.ActiveWorkbook

# Sheet 1: "Overview"
#   G2 : Latest HO Xliff Generate Date for 63564171-... .md
#        2016-08-25 13:07:17 -> 2016-08-25 13:08:15
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 13:08:15"

# Sheet 2: "zh-cn"
#   H2 : Correspond Handoff Datetime
#        2016-08-25 13:07:11 -> 2016-08-25 13:08:03
#   K2 : Correspond Handback DateTime
#        2016-08-25 13:07:35 -> 2016-08-25 13:08:29
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 13:08:03"
$wsZhCn.Range("K2").Value = "2016-08-25 13:08:29"

# Sheet 3: "de-de"
#   H2 : Correspond Handoff Datetime - shares the same text as Overview!G2,
#        so it must be updated in lockstep (2016-08-25 13:07:17 -> 2016-08-25 13:08:15)
#   K2 : Correspond Handback DateTime
#        2016-08-25 13:07:43 -> 2016-08-25 13:08:36
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 13:08:15"
$wsDeDe.Range("K2").Value = "2016-08-25 13:08:36"
